# "save works: find first value cell"
# Locate the first row whose ID column already holds "14" (the next
# record after the existing block of rows 1-13) and insert two new
# records just above it, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstValueCell = $ws.Range("A:A").Find("14")
$targetRow = $firstValueCell.Row

$insertRange = $targetRow.ToString() + ":" + ($targetRow + 1).ToString()
$ws.Rows($insertRange).Insert()

$ws.Cells.Item($targetRow, 1).Value = "111"
$ws.Cells.Item($targetRow, 2).Value = "Nod"
$ws.Cells.Item($targetRow, 3).Value = "Rod"
$ws.Cells.Item($targetRow, 4).Value = "Brazil"

$ws.Cells.Item($targetRow + 1, 1).Value = "112"
$ws.Cells.Item($targetRow + 1, 2).Value = "Mori"
$ws.Cells.Item($targetRow + 1, 3).Value = "Cox"
$ws.Cells.Item($targetRow + 1, 4).Value = "Laos"

$wb.Save()
